# Updates cryptos list data (prices and 1h volume %) to reflect the latest
# scrape, matching the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '29.904.29' }
    @{ Cell = "E2"; Value = '  -1.15%  ' }
    @{ Cell = "D3"; Value = '1.897.06' }
    @{ Cell = "E3"; Value = '  -0.74%  ' }
    @{ Cell = "E4"; Value = '  -0.12%  ' }
    @{ Cell = "D5"; Value = '''0.7558' }
    @{ Cell = "E5"; Value = '  +2.16%  ' }
    @{ Cell = "D6"; Value = '''240.10' }
    @{ Cell = "E6"; Value = '  -1.54%  ' }
    @{ Cell = "E7"; Value = '  -0.09%  ' }
    @{ Cell = "D8"; Value = '''0.3043' }
    @{ Cell = "E8"; Value = '  -2.95%  ' }
    @{ Cell = "D9"; Value = '''25.39' }
    @{ Cell = "E9"; Value = '  -6.11%  ' }
    @{ Cell = "D10"; Value = '''0.06837' }
    @{ Cell = "E10"; Value = '  -1.91%  ' }
    @{ Cell = "D11"; Value = '''0.07979' }
    @{ Cell = "E11"; Value = '  +0.10%  ' }
    @{ Cell = "B12"; Value = 'WrappedEther' }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = "D12"; Value = '1.917.07' }
    @{ Cell = "E12"; Value = '  +0.24%  ' }
    @{ Cell = "B13"; Value = 'Polygon' }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = "D13"; Value = '''0.7482' }
    @{ Cell = "E13"; Value = '  -4.27%  ' }
    @{ Cell = "D14"; Value = '''5.208' }
    @{ Cell = "E14"; Value = '  -1.62%  ' }
    @{ Cell = "D15"; Value = '''91.23' }
    @{ Cell = "E15"; Value = '  -0.52%  ' }
    @{ Cell = "D16"; Value = '29.918.27' }
    @{ Cell = "E16"; Value = '  -1.15%  ' }
    @{ Cell = "D17"; Value = '''13.95' }
    @{ Cell = "E17"; Value = '  -2.59%  ' }
    @{ Cell = "D18"; Value = '''5.959' }
    @{ Cell = "E18"; Value = '  +0.37%  ' }
    @{ Cell = "D19"; Value = '''243.10' }
    @{ Cell = "E19"; Value = '  -0.61%  ' }
    @{ Cell = "D20"; Value = '''0.000007724' }
    @{ Cell = "E20"; Value = '  -1.66%  ' }
    @{ Cell = "E21"; Value = '  -0.04%  ' }
    @{ Cell = "B22"; Value = 'BinanceUSD' }
    @{ Cell = "C22"; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' }
    @{ Cell = "D22"; Value = '''1.001' }
    @{ Cell = "E22"; Value = '  -0.09%  ' }
    @{ Cell = "B23"; Value = 'Chainlink' }
    @{ Cell = "C23"; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = "D23"; Value = '''6.952' }
    @{ Cell = "E23"; Value = '  +4.41%  ' }
    @{ Cell = "B24"; Value = 'Cosmos' }
    @{ Cell = "C24"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = "D24"; Value = '''9.235' }
    @{ Cell = "E24"; Value = '  -2.25%  ' }
    @{ Cell = "B25"; Value = 'Monero' }
    @{ Cell = "C25"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = "D25"; Value = '''165.97' }
    @{ Cell = "E25"; Value = '  +0.37%  ' }
    @{ Cell = "B26"; Value = 'EthereumClassic' }
    @{ Cell = "C26"; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = "D26"; Value = '''18.73' }
    @{ Cell = "E26"; Value = '  -1.38%  ' }
    @{ Cell = "B27"; Value = 'Stellar' }
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Cell = "D27"; Value = '''0.1304' }
    @{ Cell = "E27"; Value = '  +2.40%  ' }
    @{ Cell = "B28"; Value = 'LidoDAOToken' }
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = "D28"; Value = '''2.015' }
    @{ Cell = "E28"; Value = '  -4.53%  ' }
    @{ Cell = "B29"; Value = 'Toncoin' }
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = "D29"; Value = '''1.413' }
    @{ Cell = "E29"; Value = '  +4.22%  ' }
    @{ Cell = "B30"; Value = 'PancakeSwap' }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = "D30"; Value = '''1.519' }
    @{ Cell = "E30"; Value = '  -1.91%  ' }
    @{ Cell = "B31"; Value = 'Filecoin' }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = "D31"; Value = '''4.279' }
    @{ Cell = "E31"; Value = '  -0.95%  ' }
    @{ Cell = "B32"; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = "D32"; Value = '''4.024' }
    @{ Cell = "E32"; Value = '  -1.48%  ' }
    @{ Cell = "B33"; Value = 'Hedera' }
    @{ Cell = "C33"; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = "D33"; Value = '''0.05374' }
    @{ Cell = "E33"; Value = '  +3.53%  ' }
    @{ Cell = "B34"; Value = 'ARBITRUM' }
    @{ Cell = "C34"; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = "D34"; Value = '''1.253' }
    @{ Cell = "E34"; Value = '  -3.42%  ' }
    @{ Cell = "B35"; Value = 'ImmutableX' }
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = "D35"; Value = '''0.7257' }
    @{ Cell = "E35"; Value = '  -2.73%  ' }
    @{ Cell = "B36"; Value = 'HuobiToken' }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = "D36"; Value = '''2.719' }
    @{ Cell = "E36"; Value = '  -1.64%  ' }
    @{ Cell = "B37"; Value = 'VeChain' }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = "D37"; Value = '''0.01916' }
    @{ Cell = "E37"; Value = '  -1.44%  ' }
    @{ Cell = "B38"; Value = 'MXToken' }
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = "D38"; Value = '''2.789' }
    @{ Cell = "E38"; Value = '  -0.18%  ' }
    @{ Cell = "B39"; Value = 'FraxShare' }
    @{ Cell = "C39"; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = "D39"; Value = '''6.168' }
    @{ Cell = "E39"; Value = '  -2.94%  ' }
    @{ Cell = "B40"; Value = 'TheSandbox' }
    @{ Cell = "C40"; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = "D40"; Value = '''0.4404' }
    @{ Cell = "E40"; Value = '  -1.85%  ' }
    @{ Cell = "B41"; Value = 'Aave' }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = "D41"; Value = '''72.24' }
    @{ Cell = "E41"; Value = '  -4.48%  ' }
    @{ Cell = "B42"; Value = 'RenderToken' }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = "D42"; Value = '''1.909' }
    @{ Cell = "E42"; Value = '  -1.69%  ' }
    @{ Cell = "B43"; Value = 'PaxDollar' }
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Cell = "D43"; Value = '''1.001' }
    @{ Cell = "E43"; Value = '  +0.05%  ' }
    @{ Cell = "B44"; Value = 'TrustWalletToken' }
    @{ Cell = "C44"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = "D44"; Value = '''0.8266' }
    @{ Cell = "E44"; Value = '  -1.01%  ' }
    @{ Cell = "B45"; Value = 'Quant' }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = "D45"; Value = '''101.09' }
    @{ Cell = "E45"; Value = '  -0.09%  ' }
    @{ Cell = "B46"; Value = 'Aptos' }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = "D46"; Value = '''7.569' }
    @{ Cell = "E46"; Value = '  -1.78%  ' }
    @{ Cell = "B47"; Value = 'EnergySwap' }
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = "D47"; Value = '''9.795' }
    @{ Cell = "E47"; Value = '  -1.00%  ' }
    @{ Cell = "B48"; Value = 'RocketPoolETH' }
    @{ Cell = "C48"; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' }
    @{ Cell = "D48"; Value = '2.060.97' }
    @{ Cell = "E48"; Value = '  -0.86%  ' }
    @{ Cell = "B49"; Value = 'Elrond' }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' }
    @{ Cell = "D49"; Value = '''36.27' }
    @{ Cell = "E49"; Value = '  -3.16%  ' }
    @{ Cell = "B50"; Value = 'Cronos' }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = "D50"; Value = '''0.05978' }
    @{ Cell = "E50"; Value = '  -0.37%  ' }
    @{ Cell = "B51"; Value = 'NEARProtocol' }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = "D51"; Value = '''1.467' }
    @{ Cell = "E51"; Value = '  -0.47%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
